# Adds a "Customer Name" column to the payment-upload template sheet.
# The new column is inserted between the existing "VA Number" (C) and
# "Payment Date" (D) columns, pushing Payment Date / Journal Number /
# Payment Amount one column to the right (D:F -> E:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column at D; Excel shifts the old D:F -> E:G and the
# new column D inherits column C's formatting (font/fill/border/number
# format), matching how Excel's own "Insert Column" behaves.
$ws.Columns("D:D").Insert()

# Header (row 1) and sample value (row 2) for the new column.
$ws.Range("D1").Value = "Customer Name"
$ws.Range("D2").Value = "KLIKLELANG-Eddy susiyanto"

# Give the freshly inserted column an explicit width (it mirrors the
# "VA Number" column's width in the source workbook).
$ws.Columns("D:D").ColumnWidth = 16.33

# Move/restore the active selection to match the saved workbook state.
$ws.Range("G13").Select() | Out-Null
